$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F8").Value = 11689
$ws.Range("F26").Value = 3551
$ws.Range("F27").Value = 3551
$ws.Range("F29").Value = 798
$ws.Range("F33").Value = 982
$ws.Range("F35").Value = 65
$ws.Range("F40").Value = 2764
$ws.Range("F41").Value = 4404

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 4154
$ws.Range("F11").Value = 656

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F10").Value = 11689
$ws.Range("F21").Value = 4154
$ws.Range("F25").Value = 3551
$ws.Range("F27").Value = 798
$ws.Range("F32").Value = 982
$ws.Range("F33").Value = 65
$ws.Range("F37").Value = 4404

$wb.Save()
